# Auto-generated PowerShell COM-interop script to update the "Pais" worksheet
# to reflect the refreshed COVID-19 country data (re-sorted by total cases).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp update
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 21:22"

# Data rows: Country, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    @(4, "Estados Unidos", 356653, 19980, 19308, 326829, 8876, 900, 10516),
    @(5, "España", 135032, 3386, 40437, 81426, 6931, 528, 13169),
    @(6, "Italia", 132547, 3599, 22837, 93187, 3898, 636, 16523),
    @(7, "Alemania", 101779, 1656, 28700, 71415, 3936, 80, 1664),
    @(8, "Francia", 98010, 5171, 17250, 71849, 7072, 833, 8911),
    @(9, "China", 81708, 39, 77078, 1299, 265, 2, 3331),
    @(10, "Iran", 60500, 2274, 24236, 32525, 4083, 136, 3739),
    @(11, "Reino Unido", 51608, 3802, 135, 46100, 1559, 439, 5373),
    @(12, "Turquia", 30217, 3148, 1326, 28242, 1415, 75, 649),
    @(13, "Suiza", 21652, 552, 8056, 12834, 391, 47, 762),
    @(14, "Belgica", 20814, 1123, 3986, 15196, 1257, 185, 1632),
    @(15, "Paises Bajos", 18803, 952, 250, 16686, 1409, 101, 1867),
    @(16, "Canada", 16498, 986, 3439, 12738, 426, 41, 321),
    @(17, "Austria", 12293, 242, 3463, 8610, 250, 16, 220),
    @(18, "Portugal", 11730, 452, 140, 11279, 270, 16, 311),
    @(19, "Brasil", 11721, 467, 127, 11078, 296, 30, 516),
    @(20, "Corea del Sur", 10284, 47, 6598, 3500, 55, 3, 186),
    @(21, "Israel", 8904, 474, 585, 8262, 140, 8, 57),
    @(22, "Suecia", 7206, 376, 205, 6524, 590, 76, 477),
    @(23, "Rusia", 6343, 954, 406, 5890, 8, 2, 47),
    @(24, "Australia", 5795, 45, 2432, 3322, 96, 4, 41),
    @(25, "Noruega", 5763, 76, 32, 5655, 83, 5, 76),
    @(26, "Irlanda", 5364, 370, 25, 5165, 165, 16, 174),
    @(27, "Chile", 4815, 344, 728, 4050, 327, 3, 37),
    @(28, "India", 4778, 489, 375, 4267, 0, 18, 136),
    @(29, "Chequia", 4735, 148, 121, 4536, 84, 11, 78),
    @(30, "Dinamarca", 4681, 312, 1378, 3116, 144, 8, 187),
    @(31, "Polonia", 4413, 311, 162, 4144, 50, 13, 107),
    @(32, "Rumania", 4057, 193, 406, 3483, 179, 17, 168),
    @(33, "Malasia", 3793, 131, 1241, 2490, 102, 1, 62),
    @(34, "Pakistan", 3766, 609, 259, 3455, 17, 5, 52),
    @(35, "Ecuador", 3747, 101, 100, 3456, 156, 11, 191),
    @(36, "Filipinas", 3660, 414, 73, 3424, 1, 11, 163),
    @(37, "Japon", 3654, 0, 575, 2994, 69, 0, 85),
    @(38, "Luxemburgo", 2843, 39, 500, 2302, 33, 5, 41),
    @(39, "Arabia Saudita", 2605, 203, 551, 2016, 41, 4, 38),
    @(40, "Peru", 2561, 280, 997, 1472, 89, 9, 92),
    @(41, "Indonesia", 2491, 218, 192, 2090, 0, 11, 209),
    @(42, "Tailandia", 2220, 51, 793, 1401, 23, 3, 26),
    @(43, "Serbia", 2200, 292, 54, 2088, 101, 7, 58),
    @(44, "Finlandia", 2176, 249, 300, 1849, 81, 1, 27),
    @(45, "Mexico", 2143, 253, 633, 1416, 293, 15, 94),
    @(46, "Emiratos Arabes Unidos", 2076, 277, 167, 1898, 1, 1, 11),
    @(47, "Panama", 1988, 187, 13, 1921, 78, 8, 54),
    @(48, "Catar", 1832, 228, 131, 1697, 37, 0, 4),
    @(49, "Republica Dominicana", 1828, 83, 33, 1709, 147, 4, 86),
    @(50, "Grecia", 1755, 20, 269, 1407, 90, 6, 79),
    @(51, "Sudafrica", 1655, 0, 95, 1549, 7, 0, 11),
    @(52, "Islandia", 1562, 76, 460, 1096, 11, 2, 6),
    @(53, "Argentina", 1554, 0, 325, 1181, 94, 2, 48),
    @(54, "Colombia", 1485, 0, 88, 1362, 50, 0, 35),
    @(55, "Argelia", 1423, 103, 90, 1160, 46, 21, 173),
    @(56, "Singapur", 1375, 66, 344, 1025, 25, 0, 6),
    @(57, "Ucrania", 1319, 11, 28, 1253, 16, 1, 38),
    @(58, "Croacia", 1222, 40, 130, 1076, 36, 1, 16),
    @(59, "Egipto", 1173, 0, 247, 848, 0, 0, 78),
    @(60, "Marruecos", 1120, 99, 81, 959, 1, 10, 80),
    @(61, "Estonia", 1108, 11, 62, 1027, 14, 4, 19),
    @(62, "Nueva Zelanda", 1106, 67, 176, 929, 3, 0, 1),
    @(63, "Irak", 1031, 70, 344, 623, 0, 3, 64),
    @(64, "Eslovenia", 1021, 24, 102, 889, 30, 2, 30),
    @(65, "Moldavia", 965, 101, 37, 909, 80, 4, 19),
    @(66, "Hong Kong", 915, 24, 216, 695, 12, 0, 4),
    @(67, "Lituania", 843, 32, 8, 820, 11, 2, 15),
    @(68, "Armenia", 833, 11, 62, 763, 30, 1, 8),
    @(69, "Barein", 756, 56, 458, 294, 4, 0, 4),
    @(70, "Hungria", 744, 11, 67, 639, 17, 4, 38),
    @(71, "Crucero", 712, 0, 619, 82, 10, 0, 11),
    @(72, "Bielorrusia", 700, 138, 53, 634, 11, 5, 13),
    @(73, "Bosnia y Herzegovina", 674, 20, 47, 598, 4, 6, 29),
    @(74, "Kuwait", 665, 109, 103, 561, 20, 0, 1),
    @(75, "Kazajistan", 662, 78, 46, 610, 16, 0, 6),
    @(76, "Camerun", 658, 8, 17, 632, 0, 0, 9),
    @(77, "Azerbaiyan", 641, 57, 44, 590, 11, 0, 7),
    @(78, "Tunez", 596, 22, 5, 569, 39, 0, 22),
    @(79, "Republica de Macedonia", 570, 15, 30, 519, 15, 3, 21),
    @(80, "Bulgaria", 549, 18, 39, 488, 26, 2, 22),
    @(81, "Letonia", 542, 9, 16, 525, 5, 0, 1),
    @(82, "Libano", 541, 14, 60, 462, 27, 1, 19),
    @(83, "Eslovaquia", 534, 49, 8, 524, 3, 1, 2),
    @(84, "Principado de Andorra", 501, 0, 26, 457, 12, 0, 18),
    @(85, "Costa Rica", 467, 13, 18, 447, 14, 0, 2),
    @(86, "Republica de Chipre", 465, 19, 45, 411, 11, 0, 9),
    @(87, "Uzbekistan", 457, 115, 30, 425, 8, 0, 2),
    @(88, "Uruguay", 406, 0, 104, 296, 14, 0, 6),
    @(89, "Albania", 377, 16, 116, 240, 7, 1, 21),
    @(90, "Taiwan", 373, 10, 57, 311, 0, 0, 5),
    @(91, "Afganistan", 367, 18, 18, 342, 0, 0, 7),
    @(92, "Cuba", 350, 30, 18, 323, 11, 1, 9),
    @(93, "Reunion", 349, 5, 40, 309, 4, 0, 0),
    @(94, "Jordania", 349, 4, 126, 217, 5, 1, 6),
    @(95, "Burkina Faso", 345, 0, 90, 238, 0, 0, 17),
    @(96, "Oman", 331, 33, 61, 268, 3, 0, 2),
    @(97, "Honduras", 298, 30, 6, 270, 10, 0, 22),
    @(98, "San Marino", 266, 0, 35, 199, 14, 0, 32),
    @(99, "Costa de Marfil", 261, 0, 37, 221, 0, 0, 3),
    @(100, "Estado de Palestina", 253, 16, 25, 227, 0, 0, 1),
    @(101, "Vietnam", 245, 4, 95, 150, 8, 0, 0),
    @(102, "Mauricio", 244, 17, 7, 230, 2, 0, 7),
    @(103, "Malta", 241, 14, 5, 236, 3, 0, 0),
    @(104, "Montenegro", 233, 19, 1, 230, 4, 0, 2),
    @(105, "Nigeria", 232, 0, 33, 194, 2, 0, 5),
    @(106, "Senegal", 226, 4, 92, 132, 1, 0, 2),
    @(107, "Kirguistan", 216, 69, 33, 179, 5, 3, 4),
    @(108, "Ghana", 214, 0, 31, 178, 2, 0, 5),
    @(109, "Georgia", 188, 14, 39, 147, 6, 0, 2),
    @(110, "Niger", 184, 0, 13, 161, 0, 0, 10),
    @(111, "Bolivia", 183, 26, 2, 170, 3, 1, 11),
    @(112, "Islas Feroe", 183, 2, 107, 76, 1, 0, 0),
    @(113, "Sri Lanka", 178, 2, 38, 135, 5, 0, 5),
    @(114, "Consejo Danes para los Refugiados", 161, 7, 5, 138, 0, 0, 18),
    @(115, "Venezuela", 159, 0, 52, 100, 6, 0, 7),
    @(116, "Kenia", 158, 16, 4, 148, 2, 2, 6),
    @(117, "Martinica", 149, 0, 50, 95, 21, 0, 4),
    @(118, "Mayotte", 147, 0, 14, 131, 3, 0, 2),
    @(119, "Isla de Man", 139, 12, 55, 83, 0, 0, 1),
    @(120, "Guadalupe", 135, 0, 31, 97, 14, 0, 7),
    @(121, "Brunei", 135, 0, 82, 52, 3, 0, 1),
    @(122, "Banglades", 123, 35, 33, 78, 1, 3, 12),
    @(123, "Guinea", 121, 0, 5, 116, 0, 0, 0),
    @(124, "Camboya", 114, 0, 53, 61, 1, 0, 0),
    @(125, "Paraguay", 113, 9, 12, 96, 8, 2, 5),
    @(126, "Gibraltar", 109, 6, 52, 57, 0, 0, 0),
    @(127, "Ruanda", 105, 1, 4, 101, 0, 0, 0),
    @(128, "Trinidad yTobago", 105, 1, 1, 96, 0, 1, 8),
    @(129, "Monaco", 77, 4, 4, 72, 4, 0, 1),
    @(130, "Liechtenstein", 77, 0, 55, 21, 0, 0, 1),
    @(131, "Madagascar", 72, 0, 2, 70, 6, 0, 0),
    @(132, "Aruba", 71, 7, 2, 69, 0, 0, 0),
    @(133, "Guatemala", 70, 0, 15, 52, 3, 0, 3),
    @(134, "El Salvador", 69, 7, 5, 60, 4, 1, 4),
    @(135, "Guayana Francesa", 68, 0, 27, 41, 1, 0, 0),
    @(136, "Republica de Yibuti", 59, 0, 9, 50, 0, 0, 0),
    @(137, "Jamaica", 58, 0, 8, 47, 0, 0, 3),
    @(138, "Barbados", 56, 0, 6, 49, 4, 0, 1),
    @(139, "Uganda", 52, 0, 0, 52, 0, 0, 0),
    @(140, "Togo", 52, 8, 22, 27, 0, 0, 3),
    @(141, "Mali", 47, 2, 1, 41, 0, 0, 5),
    @(142, "Congo", 45, 0, 2, 38, 0, 0, 5),
    @(143, "Etiopia", 44, 1, 4, 38, 1, 0, 2),
    @(144, "Macao", 44, 0, 10, 34, 0, 0, 0),
    @(145, "Polinesia Francesa", 42, 1, 0, 42, 0, 0, 0),
    @(146, "Islas Caimanes", 39, 0, 1, 37, 0, 0, 1),
    @(147, "Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @(148, "Zambia", 39, 0, 5, 33, 0, 0, 1),
    @(149, "San Martin (Parte Holandesa)", 37, 12, 1, 30, 0, 2, 6),
    @(150, "Bermudas", 37, 0, 14, 23, 0, 0, 0),
    @(151, "Guam", 32, 0, 0, 31, 0, 0, 1),
    @(152, "San Martin (Parte Francesa)", 32, 0, 7, 23, 6, 0, 2),
    @(153, "Eritrea", 29, 0, 0, 29, 0, 0, 0),
    @(154, "Guyana", 29, 5, 0, 25, 0, 0, 4),
    @(155, "Bahamas", 29, 0, 4, 20, 1, 0, 5),
    @(156, "Haiti", 24, 3, 0, 23, 0, 0, 1),
    @(157, "Gabon", 24, 3, 1, 22, 0, 0, 1),
    @(158, "Tanzania", 24, 2, 3, 20, 0, 0, 1),
    @(159, "Benin", 23, 1, 5, 17, 0, 1, 1),
    @(160, "Birmania", 22, 1, 0, 21, 0, 0, 1),
    @(161, "Siria", 19, 0, 2, 15, 0, 0, 2),
    @(162, "Maldivas", 19, 0, 13, 6, 0, 0, 0),
    @(163, "Guinea-Bisau", 18, 0, 0, 18, 0, 0, 0),
    @(164, "Nueva Caledonia", 18, 0, 1, 17, 0, 0, 0),
    @(165, "Libia", 18, 0, 0, 17, 0, 0, 1),
    @(166, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @(167, "Guinea Ecuatorial", 16, 0, 3, 13, 0, 0, 0),
    @(168, "Namibia", 16, 0, 3, 13, 0, 0, 0),
    @(169, "Angola", 16, 2, 2, 12, 0, 0, 2),
    @(170, "Antigua y Barbuda", 15, 0, 0, 15, 1, 0, 0),
    @(171, "Mongolia", 15, 1, 2, 13, 0, 0, 0),
    @(172, "Fiyi", 14, 2, 0, 14, 0, 0, 0),
    @(173, "Dominica", 14, 0, 0, 14, 0, 0, 0),
    @(174, "Santa Lucia", 14, 0, 1, 13, 0, 0, 0),
    @(175, "Liberia", 14, 1, 3, 8, 0, 0, 3),
    @(176, "Curazao", 13, 2, 5, 7, 0, 0, 1),
    @(177, "Granada", 12, 0, 0, 12, 2, 0, 0),
    @(178, "Laos", 12, 1, 0, 12, 0, 0, 0),
    @(179, "Sudan", 12, 0, 2, 8, 0, 0, 2),
    @(180, "Seychelles", 11, 1, 0, 11, 0, 0, 0),
    @(181, "Groenlandia", 11, 0, 4, 7, 0, 0, 0),
    @(182, "San Cristobal y Nieves", 10, 0, 0, 10, 0, 0, 0),
    @(183, "Mozambique", 10, 0, 1, 9, 0, 0, 0),
    @(184, "Surinam", 10, 0, 0, 9, 0, 0, 1),
    @(185, "Suazilandia", 10, 1, 4, 6, 0, 0, 0),
    @(186, "Republica del Chad", 9, 0, 0, 9, 0, 0, 0),
    @(187, "Nepal", 9, 0, 1, 8, 0, 0, 0),
    @(188, "Zimbabue", 9, 0, 0, 8, 0, 0, 1),
    @(189, "Montserrat", 9, 0, 0, 7, 0, 0, 2),
    @(190, "Republica de Africa Central", 8, 0, 0, 8, 0, 0, 0),
    @(191, "Santa Sede", 7, 0, 0, 7, 0, 0, 0),
    @(192, "Belice", 7, 2, 0, 6, 1, 1, 1),
    @(193, "San Vicente y las Granadinas", 7, 0, 1, 6, 0, 0, 0),
    @(194, "Somalia", 7, 0, 1, 6, 0, 0, 0),
    @(195, "Cabo Verde", 7, 0, 1, 5, 0, 0, 1),
    @(196, "Sierra Leona", 6, 0, 0, 6, 0, 0, 0),
    @(197, "Nicaragua", 6, 0, 0, 5, 0, 0, 1),
    @(198, "Botsuana", 6, 0, 0, 5, 0, 0, 1),
    @(199, "San Bartolome", 6, 0, 1, 5, 0, 0, 0),
    @(200, "Mauritania", 6, 0, 2, 3, 0, 0, 1),
    @(201, "Malaui", 5, 1, 0, 5, 1, 0, 0),
    @(202, "Islas Turcas y Caicos", 5, 0, 0, 4, 0, 0, 1),
    @(203, "Butan", 5, 0, 2, 3, 0, 0, 0),
    @(204, "Santo Tome y Principe", 4, 4, 0, 4, 0, 0, 0),
    @(205, "Sahara Occidental", 4, 0, 0, 4, 0, 0, 0),
    @(206, "Gambia", 4, 0, 2, 1, 0, 0, 1),
    @(207, "Burundi", 3, 0, 0, 3, 0, 0, 0),
    @(208, "Anguila", 3, 0, 0, 3, 0, 0, 0),
    @(209, "Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0),
    @(210, "Papua Nueva Guinea", 2, 1, 0, 2, 0, 0, 0),
    @(211, "Bonaire, San Eustaquio y Saba", 2, 0, 0, 2, 0, 0, 0),
    @(212, "Islas Malvinas", 2, 0, 0, 2, 0, 0, 0),
    @(213, "Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
    @(214, "Sudan del Sur", 1, 0, 0, 1, 0, 0, 0),
    @(215, "San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = [double]$row[2]
    $ws.Cells.Item($r, 3).Value = [double]$row[3]
    $ws.Cells.Item($r, 4).Value = [double]$row[4]
    $ws.Cells.Item($r, 5).Value = [double]$row[5]
    $ws.Cells.Item($r, 6).Value = [double]$row[6]
    $ws.Cells.Item($r, 7).Value = [double]$row[7]
    $ws.Cells.Item($r, 8).Value = [double]$row[8]
}

Write-Output "Done: updated $($data.Count) rows"
